$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, styled like the other header cells (copy style from E1)
$ws.Range("F1").Value = "Points"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Point name values for rows 2-13 (A through L)
$points = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L")
for ($i = 0; $i -lt $points.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $points[$i]
}
